# Auto-generated edit script applying cryptos.xlsx row updates from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.674.80'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '3.527.57'
$ws.Range('E3').Value = '  -0.72%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''614.43'
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('D6').Value = '''173.94'
$ws.Range('E6').Value = '  +1.12%  '
$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').Value = '3.520.95'
$ws.Range('E7').Value = '  -0.84%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = '''0.610'
$ws.Range('E8').Value = '  -1.36%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('E10').Value = '  +0.08%  '
$ws.Range('D11').Value = '''7.23'
$ws.Range('E11').Value = '  +4.26%  '
$ws.Range('D12').Value = '''0.589'
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('D13').Value = '''46.57'
$ws.Range('E13').Value = '  -0.34%  '
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('D15').Value = '4.097.72'
$ws.Range('E15').Value = '  -0.64%  '
$ws.Range('D16').Value = '''8.46'
$ws.Range('E16').Value = '  +0.83%  '
$ws.Range('D17').Value = '''616.63'
$ws.Range('E17').Value = '  -0.29%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.531.11'
$ws.Range('E18').Value = '  -0.52%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '70.675.77'
$ws.Range('E19').Value = '  -0.29%  '
$ws.Range('E20').Value = '  +1.71%  '
$ws.Range('D21').Value = '''17.78'
$ws.Range('E21').Value = '  +2.18%  '
$ws.Range('D22').Value = '''0.885'
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('D23').Value = '''9.00'
$ws.Range('E23').Value = '  -5.34%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').Value = '''15.76'
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '''98.70'
$ws.Range('E25').Value = '  +1.99%  '
$ws.Range('E26').Value = '  -0.67%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').Value = '''2.61'
$ws.Range('E28').Value = '  -0.14%  '
$ws.Range('D29').Value = '''33.94'
$ws.Range('E29').Value = '  +1.19%  '
$ws.Range('D30').Value = '''9.20'
$ws.Range('E30').Value = '  +1.45%  '
$ws.Range('B31').Value = 'Stacks'
$ws.Range('C31').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D31').Value = '''3.03'
$ws.Range('E31').Value = '  -1.85%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '''8.16'
$ws.Range('E32').Value = '  -3.80%  '
$ws.Range('D33').Value = '''1.31'
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('D34').Value = '''6.86'
$ws.Range('E34').Value = '  -2.17%  '
$ws.Range('D35').Value = '''620.90'
$ws.Range('E35').Value = '  +8.19%  '
$ws.Range('D36').Value = '''0.100'
$ws.Range('E36').Value = '  -0.69%  '
$ws.Range('D37').Value = '''0.0492'
$ws.Range('E37').Value = '  +4.34%  '
$ws.Range('E38').Value = '  +0.37%  '
$ws.Range('D39').Value = '''3.50'
$ws.Range('E39').Value = '  -3.14%  '
$ws.Range('D40').Value = '''57.02'
$ws.Range('E40').Value = '  -1.15%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').Value = '  +1.77%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '3.378.05'
$ws.Range('E43').Value = '  +0.66%  '
$ws.Range('B44').Value = 'PEPE'
$ws.Range('C44').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D44').Value = '0.0₃0739'
$ws.Range('E44').Value = '  +4.72%  '
$ws.Range('E45').Value = '  -2.58%  '
$ws.Range('D46').Value = '''2.93'
$ws.Range('E46').Value = '  -2.08%  '
$ws.Range('D47').Value = '''32.29'
$ws.Range('E47').Value = '  -2.24%  '
$ws.Range('D48').Value = '''2.57'
$ws.Range('E48').Value = '  -2.21%  '
$ws.Range('E49').Value = '  +0.52%  '
$ws.Range('D50').Value = '''133.82'
$ws.Range('E50').Value = '  +0.02%  '
